# Adds a new weekly record at row 5 (pushing all existing data rows down
# by one) for "Terminal Hortofrutícola Agro Chillán - Haba".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 5, shifting rows 5:33 down
# to 6:34 (this also extends the sheet dimension from R33 to R34 and
# carries the existing D-column date style down with the shifted cells).
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly entry.
$ws.Cells.Item(5,1).Value = 7
$ws.Cells.Item(5,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(5,3).Value = "Ñuble"
$ws.Cells.Item(5,4).Value = "11/30/2021"
$ws.Cells.Item(5,5).Value = 16
$ws.Cells.Item(5,6).Value = 100112026
$ws.Cells.Item(5,7).Value = "Haba"
$ws.Cells.Item(5,8).Value = "Sin especificar"
$ws.Cells.Item(5,9).Value = "Primera"
$ws.Cells.Item(5,10).Value = 100
$ws.Cells.Item(5,11).Value = 6000
$ws.Cells.Item(5,12).Value = 7000
$ws.Cells.Item(5,13).Value = 6500
$ws.Cells.Item(5,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(5,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(5,16).Value = 260
$ws.Cells.Item(5,17).Value = 25
$ws.Cells.Item(5,18).Value = "Hortaliza"
